$d = $word.ActiveDocument

# Replace "(1 + 2/ 3/4/ 5) " with "(1 + 2/3/4/5) " everywhere in the document.
$d.Content.Find.Execute("2/ 3/4/ 5", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2/3/4/5", 2)
